$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sprint 2 backlog review edits ---

# Row 8: clarify the outgoing-orders backlog item to also cover incoming orders
$ws.Range("C8").Value = "be able to view a detailed list of orders that are outgoing/incoming"

# Row 9: replace the old "track when things are received and shipped" item with a
# new "view history a single item" item, and reuse the "we know where rentals are"
# justification that used to sit on the (now removed) row.
$ws.Range("B9").Value = "Manager"
$ws.Range("C9").Value = "view history a single item"
$ws.Range("D9").Value = " we know where rentals are "

# Row 10: the old row-11 "manage the stock" pairing shifts up and is re-prioritized
$ws.Range("A10").Value = 2
$ws.Range("B10").Value = "Manager"
$ws.Range("C10").Value = "manage the stock"
$ws.Range("D10").Value = " we always have a good inventory"

# Row 11: the old row-10 "view records of items..." pairing shifts up and is re-prioritized
$ws.Range("A11").Value = 2
$ws.Range("B11").Value = "Manager"
$ws.Range("C11").Value = "view records of items tracked by which employee"
$ws.Range("D11").Value = " I can find items"

# --- Re-apply the cell formatting so styles land on the same cells they did before ---
# B9 picks up the "theme colour" look that D9 already had
$ws.Range("D9").Copy() | Out-Null
$ws.Range("B9").PasteSpecial(-4122) | Out-Null

# D10 picks up that same "theme colour" look
$ws.Range("D9").Copy() | Out-Null
$ws.Range("D10").PasteSpecial(-4122) | Out-Null

# B11 and D11 drop back to the plain look that B10/C10 already use
$ws.Range("B10").Copy() | Out-Null
$ws.Range("B11").PasteSpecial(-4122) | Out-Null
$ws.Range("D11").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- Column C is now wider to fit the longer item-8 text ---
$ws.Columns("C:C").AutoFit() | Out-Null

# --- Re-run the same Data > Sort the sheet already used, now over the trimmed A2:D23 range ---
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A2:A23"))
$ws.Sort.SetRange($ws.Range("A2:D23"))
$ws.Sort.Header = 0
$ws.Sort.Apply()

# --- Restore the cursor to where the reviewer left it ---
$ws.Range("B6").Select()
